# Applies the "updated PvsI with model fitting" edit: recomputed volume
# (T), rate.abs (Z), rate.a.spec (AB) and rate.output (AD) values for the
# 50_rates respirometry sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new cell values (only cells that actually change per the diff)
$updates = @{
    2  = @{ T = 0.1405756097560976;  Z = -0.0728187638146726;  AB = -293.8047393614331;  AD = -293.8047393614331 }
    3  = @{ T = 0.144009756097561;   Z = -0.1143681392315392;  AB = -753.4296980080328;  AD = -753.4296980080328 }
    4  = @{ T = 0.1462634146341464;  Z = -0.09226426848590945; AB = -591.05214797052;     AD = -591.05214797052 }
    5  = @{ T = 0.1500390243902439;  Z = -0.07314360166293993; AB = -456.8013575459984;   AD = -456.8013575459984 }
    6  = @{ T = 0.1463707317073171;  Z = -0.09728375923564792; AB = -607.7298285511221;  AD = -607.7298285511221 }
    7  = @{ T = 0.1465658536585366;  Z = -0.1122908035996341;  AB = -560.321756497604;   AD = -560.321756497604 }
    8  = @{ T = 0.1449658536585366;  Z = -0.08609147818295118; AB = -347.2643588830868;  AD = -347.2643588830868 }
    9  = @{ T = 0.1544 }
    10 = @{ T = 0.1405756097560976;  Z = 0.05890205722984346;  AB = 237.6544542325692;   AD = 237.6544542325692 }
    11 = @{ T = 0.144009756097561;   Z = 0.05732404840868841;  AB = 377.6369955072732;   AD = 377.6369955072732 }
    12 = @{ T = 0.1462634146341464;  Z = 0.03715723579856357;  AB = 238.0321698940484;   AD = 238.0321698940484 }
    13 = @{ T = 0.1500390243902439;  Z = 0.07295919816265506;  AB = 455.6497083606802;   AD = 455.6497083606802 }
    14 = @{ T = 0.1463707317073171;  Z = -0.01163180857724525; AB = -72.66369112305458;  AD = -72.66369112305458 }
    15 = @{ T = 0.1465658536585366;  Z = 0.08360906481748334;  AB = 417.2022690717208;   AD = 417.2022690717208 }
    16 = @{ T = 0.1449658536585366;  Z = 0.09039440703388167;  AB = 364.6209412101147;   AD = 364.6209412101147 }
    17 = @{ T = 0.1544;              Z = -0.0005550390967598329 }
}

foreach ($row in $updates.Keys) {
    $cellset = $updates[$row]
    foreach ($col in $cellset.Keys) {
        $ws.Range("$col$row").Value = $cellset[$col]
    }
}
